$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.794.70"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "2.353.79"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.668"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.79%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.19%  "

$ws.Range("E12").Value = "  +8.54%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("E15").Value = "  -2.72%  "

$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").Value = "2.346.52"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").Value = "43.756.33"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "252.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +3.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.24%  "

$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "

$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  -2.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.129"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0747"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.49%  "

$ws.Range("E34").Value = "  -4.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +3.96%  "

$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("E39").Value = "  -1.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.26%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.57%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.203"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("E44").Value = "  -6.16%  "

$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("E51").Value = "  +1.99%  "
